$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style from H1 (bold, bordered, centered header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-22
$iValues = @(8, 6, 5, 8, 9, 9, 6, 7, 4, 5, 9, 11, 6, 5, 3, 3, 5, 5, 4, 7, 4)
$jValues = @(8, 7, 5, 8, 9, 9, 7, 7, 4, 6, 9, 11, 7, 5, 5, 4, 6, 8, 5, 8, 6)

for ($r = 0; $r -lt 21; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
